$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data (rows 1-23) shifts down to rows 2-24.
$ws.Rows.Item(1).Insert()

# Set header values
$ws.Range("A1").Value = "iddepartamento"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "idsuperior"

# Match the header row style: right-aligned text
$ws.Range("A1:C1").HorizontalAlignment = -4152  # xlRight

# Column B width
$ws.Columns.Item(2).ColumnWidth = 45.85546875

# View adjustments: scroll so row 10 is at the top of the visible window, then
# select D9 to match the saved selection state.
try {
    $excel.ActiveWindow.ScrollRow = 10
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("D9").Select() | Out-Null
